# Update "想去人数" (interested count) values in the F column of the
# "展览" and "全部类型" sheets to reflect newly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 366
    $ws.Range("F3").Value = 242
    $ws.Range("F4").Value = 81
    $ws.Range("F5").Value = 297
}
